$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 2266.5
$ws.Range("I34").Value = 2266.5
$ws.Range("K34").Value = 2266.5
$ws.Range("M34").Value = -2063.5
# Row 36
$ws.Range("H36").Value = 2266.5
$ws.Range("I36").Value = 2266.5
$ws.Range("K36").Value = 2266.5
$ws.Range("M36").Value = -1551.5
# Row 69
$ws.Range("H69").Value = 6160.619
$ws.Range("J69").Value = 6742.222
$ws.Range("L69").Value = 20226.666
$ws.Range("N69").Value = -21974.666
# Row 72
$ws.Range("H72").Value = 6160.619
$ws.Range("J72").Value = 6742.222
$ws.Range("L72").Value = 60679.998
$ws.Range("N72").Value = -69415.99799999999
# Row 80
$ws.Range("H80").Value = 1759.8
$ws.Range("J80").Value = 1666.3334
$ws.Range("L80").Value = 4999.0002
$ws.Range("I80").Value = 1900
$ws.Range("K80").Value = 5700
$ws.Range("M80").Value = -4702
$ws.Range("N80").Value = -6995.0002
# Row 83
$ws.Range("H83").Value = 1759.8
$ws.Range("J83").Value = 1666.3334
$ws.Range("L83").Value = 14997.0006
$ws.Range("I83").Value = 1900
$ws.Range("K83").Value = 17100
$ws.Range("M83").Value = -12108
$ws.Range("N83").Value = -24981.0006
# Row 88
$ws.Range("H88").Value = 1529.4546
$ws.Range("J88").Value = 1904.8572
$ws.Range("L88").Value = 1904.8572
$ws.Range("N88").Value = -2716.8572
# Row 91
$ws.Range("H91").Value = 1529.4546
$ws.Range("J91").Value = 1904.8572
$ws.Range("L91").Value = 1904.8572
$ws.Range("N91").Value = -4712.8572
# Row 111
$ws.Range("H111").Value = 1582.9231
$ws.Range("I111").Value = 461.4
$ws.Range("K111").Value = 1384.2
$ws.Range("M111").Value = 1682.8
# Row 112
$ws.Range("H112").Value = 2816.6667
$ws.Range("J112").Value = 2975
$ws.Range("L112").Value = 8925
$ws.Range("N112").Value = -11141
# Row 125
$ws.Range("H125").Value = 955.8333
$ws.Range("J125").Value = 735
$ws.Range("L125").Value = 6615
$ws.Range("N125").Value = -11535
# Row 134
$ws.Range("H134").Value = 205000
$ws.Range("J134").Value = 205000
$ws.Range("L134").Value = 205000
$ws.Range("N134").Value = -215140
# Row 137
$ws.Range("H137").Value = 2516.547
$ws.Range("I137").Value = 1583.7587
$ws.Range("K137").Value = 4751.2761
$ws.Range("M137").Value = -2201.2761
# Row 138
$ws.Range("H138").Value = 3239.75
$ws.Range("J138").Value = 3198.3333
$ws.Range("L138").Value = 9594.999899999999
$ws.Range("I138").Value = 3281.1667
$ws.Range("K138").Value = 9843.500100000001
$ws.Range("M138").Value = -4703.500100000001
$ws.Range("N138").Value = -19874.9999

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3736.375
$ws.Range("I61").Value = 3736.375
$ws.Range("K61").Value = 3736.375
$ws.Range("M61").Value = -3524.375
# Row 136
$ws.Range("H136").Value = 3736.375
$ws.Range("I136").Value = 3736.375
$ws.Range("K136").Value = 11209.125
$ws.Range("M136").Value = -8659.125

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 603.6667
$ws.Range("J16").Value = 552
$ws.Range("L16").Value = 552
$ws.Range("I16").Value = 645
$ws.Range("K16").Value = 645
$ws.Range("M16").Value = -358
$ws.Range("N16").Value = -1126
# Row 31
$ws.Range("H31").Value = 5757.5957
$ws.Range("J31").Value = 6020.1396
$ws.Range("L31").Value = 6020.1396
$ws.Range("I31").Value = 2935.25
$ws.Range("K31").Value = 2935.25
$ws.Range("M31").Value = -2640.25
$ws.Range("N31").Value = -6610.1396
# Row 32
$ws.Range("H32").Value = 2853.3333
$ws.Range("J32").Value = 3050
$ws.Range("L32").Value = 3050
$ws.Range("I32").Value = 2755
$ws.Range("K32").Value = 2755
$ws.Range("M32").Value = -2439
$ws.Range("N32").Value = -3682
# Row 34
$ws.Range("H34").Value = 5757.5957
$ws.Range("J34").Value = 6020.1396
$ws.Range("L34").Value = 6020.1396
$ws.Range("I34").Value = 2935.25
$ws.Range("K34").Value = 2935.25
$ws.Range("M34").Value = -2733.25
$ws.Range("N34").Value = -6424.1396
# Row 58
$ws.Range("H58").Value = 3104.4666
$ws.Range("I58").Value = 2504
$ws.Range("K58").Value = 2504
$ws.Range("M58").Value = -2301
# Row 113
$ws.Range("H113").Value = 603.6667
$ws.Range("J113").Value = 552
$ws.Range("L113").Value = 552
$ws.Range("I113").Value = 645
$ws.Range("K113").Value = 645
$ws.Range("M113").Value = 1525
$ws.Range("N113").Value = -4892
# Row 134
$ws.Range("H134").Value = 1624.375
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("I134").Value = 1624.375
$ws.Range("K134").Value = 4873.125
$ws.Range("M134").Value = -2338.125
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 3104.4666
$ws.Range("I136").Value = 2504
$ws.Range("K136").Value = 7512
$ws.Range("M136").Value = -4962
# Row 141
$ws.Range("H141").Value = 51130.4
$ws.Range("J141").Value = 51130.4
$ws.Range("L141").Value = 51130.4
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -61490.4

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 819.4737
$ws.Range("J5").Value = 779.4091
$ws.Range("L5").Value = 2338.2273
$ws.Range("I5").Value = 874.5625
$ws.Range("K5").Value = 2623.6875
$ws.Range("M5").Value = -2511.6875
$ws.Range("N5").Value = -2562.2273
# Row 12
$ws.Range("H12").Value = 177.58333
$ws.Range("J12").Value = 195.5
$ws.Range("L12").Value = 586.5
$ws.Range("N12").Value = -932.5
# Row 74
$ws.Range("H74").Value = 8757.5
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 8757.5
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 129
$ws.Range("H129").Value = 2074.3333
$ws.Range("J129").Value = 2374.5
$ws.Range("L129").Value = 7123.5
$ws.Range("I129").Value = 1474
$ws.Range("K129").Value = 4422
$ws.Range("M129").Value = 578
$ws.Range("N129").Value = -17123.5
# Row 135
$ws.Range("H135").Value = 819.4737
$ws.Range("J135").Value = 779.4091
$ws.Range("L135").Value = 7014.6819
$ws.Range("I135").Value = 874.5625
$ws.Range("K135").Value = 7871.0625
$ws.Range("M135").Value = -5336.0625
$ws.Range("N135").Value = -12084.6819

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2048.0715
$ws.Range("I102").Value = 1864.4166
$ws.Range("K102").Value = 1864.4166
$ws.Range("M102").Value = -242.4166
# Row 113
$ws.Range("H113").Value = 4722
$ws.Range("I113").Value = 2416.3333
$ws.Range("K113").Value = 2416.3333
$ws.Range("M113").Value = -246.3332999999998
# Row 132
$ws.Range("H132").Value = 8999
$ws.Range("J132").Value = 8999
$ws.Range("L132").Value = 26997
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -32057

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1151.2
$ws.Range("I16").Value = 1151.2
$ws.Range("K16").Value = 1151.2
$ws.Range("M16").Value = -981.2
# Row 31
$ws.Range("H31").Value = 320.8
$ws.Range("J31").Value = 499
$ws.Range("L31").Value = 499
$ws.Range("I31").Value = 276.25
$ws.Range("K31").Value = 276.25
$ws.Range("M31").Value = -28.25
$ws.Range("N31").Value = -995
# Row 93
$ws.Range("H93").Value = 835.75
$ws.Range("I93").Value = 781
$ws.Range("K93").Value = 781
$ws.Range("M93").Value = 467
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 22576.584
$ws.Range("J41").Value = 21773.75
$ws.Range("L41").Value = 21773.75
$ws.Range("N41").Value = -22553.75
# Row 49
$ws.Range("H49").Value = 15062
$ws.Range("J49").Value = 15062
$ws.Range("L49").Value = 15062
$ws.Range("N49").Value = -15522
# Row 81
$ws.Range("H81").Value = 6648.75
$ws.Range("I81").Value = 3031.6667
$ws.Range("K81").Value = 6063.3334
$ws.Range("M81").Value = -5002.3334
# Row 84
$ws.Range("H84").Value = 6648.75
$ws.Range("I84").Value = 3031.6667
$ws.Range("K84").Value = 30316.667
$ws.Range("M84").Value = -25012.667
# Row 126
$ws.Range("H126").Value = 6430.615
$ws.Range("I126").Value = 3649.5
$ws.Range("K126").Value = 10948.5
$ws.Range("M126").Value = -8478.5
# Row 136
$ws.Range("H136").Value = 2901.3845
$ws.Range("I136").Value = 2401.7144
$ws.Range("K136").Value = 7205.1432
$ws.Range("M136").Value = -4655.1432
# Row 141
$ws.Range("H141").Value = 229998
